# Textbox response formatting fix
# Renames task-order sheets and refreshes the stimulus filenames listed in column B.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16511686580825233"
$ws1.Range("B2").Value = "go_stims-1651168658052636.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686580669339.csv"
$ws1.Range("B4").Value = "go_stims-16511686580699298.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686580815232.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1651168660004697"
$ws2.Range("B2").Value = "ZB-match_2-16511686586308818.csv"
$ws2.Range("B3").Value = "ZB-match_4-16511686583058248.csv"
$ws2.Range("B4").Value = "OB-16511686590193365.csv"
$ws2.Range("B5").Value = "OB-16511686586707425.csv"
$ws2.Range("B6").Value = "TB-16511686599789274.csv"
$ws2.Range("B7").Value = "TB-1651168659832636.csv"
$ws2.Range("B8").Value = "ZB-match_9-16511686584228742.csv"
$ws2.Range("B9").Value = "TB-1651168659288478.csv"
$ws2.Range("B10").Value = "OB-16511686591856558.csv"

# --- Sheet 3: RS_TO (name change only) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16511686600056565"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1651168660067454"
$ws4.Range("B2").Value = "MM_stims-16511686600200593.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686600066543.csv"
$ws4.Range("B4").Value = "MM_stims-16511686600507755.csv"
$ws4.Range("B5").Value = "ZM_stims-16511686600200593.csv"
$ws4.Range("B6").Value = "MM_stims-1651168660066484.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686600507755.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511686601556132"
$ws5.Range("B2").Value = "vSAT_stims-16511686601137147.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511686601407137.csv"
$ws5.Range("B4").Value = "SAT_stims-1651168660097286.csv"
$ws5.Range("B5").Value = "SAT_stims-1651168660071232.csv"
